# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.97 = 7179.11 pesos`n✅ 7179.11 pesos = 1.97 = 930.53 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 508
$wsTasas.Range("O10").Value = 3646.99
$wsTasas.Range("N12").Value = 3650
$wsTasas.Range("O12").Value = 473.1
